$d = $word.ActiveDocument

# 1. Remove the leftover "_GoBack" bookmark (empty bookmark in an otherwise
#    empty paragraph near the end of section 2.1 / before "1 Einleitung").
$d.Bookmarks("_GoBack").Delete()

# 2. Append a new sentence, as its own run, right after the existing
#    "Die Entwicklung eines Programmes ... Klingebiel. " run. We deliberately
#    avoid simply concatenating into the existing run's text (that would
#    merge both sentences into a single <w:r>) by instead splitting the
#    paragraph in two, typing the new sentence into the freshly created
#    paragraph, and then re-joining the two paragraphs by deleting the
#    paragraph mark between them. That mirrors how Word keeps the original
#    run untouched and places newly typed text into a brand-new run.
$anchor = $d.Content
$anchor.Find.Execute("Die Entwicklung eines Programmes gehört zur Aufgabe des Informatik Kurses von Herrn Klingebiel. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$anchor.Collapse(0)
$anchor.Move(4, 1) | Out-Null
$anchor.InsertAfter("Wir haben uns dabei für ein Spiel entschieden")

$joiner = $d.Content
$joiner.Find.Execute("Die Entwicklung eines Programmes gehört zur Aufgabe des Informatik Kurses von Herrn Klingebiel. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$joiner.Collapse(0)
$joiner.MoveEnd(4, 1) | Out-Null
$joiner.Delete()
